# "First draft for CIVE 596"
# Update the reach used for the width/slope computations on the N2 sheet
# (top-right summary table) from stations B8/B21 to stations B12/B17, fix
# up the dependent formula for average thalweg slope to match, give the
# slope ratio cell (K3) two decimal places, and leave the selection on G5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("N2")
$ws.Activate()

# G3: total reach length now measured between rows 12 and 17 instead of 8 and 21
$ws.Range("G3").Formula = "=ABS(B12)+ABS(B17)"

# G4: matching change to the average elevation-drop formula (rows 12 & 17)
$ws.Range("G4").Formula = "=(-((E17*12)+F17)+((E12*12)+F12))/12"

# K3 (H4/H3 slope ratio) gets a new "0.00" number format
$ws.Range("K3").NumberFormat = "0.00"

# Leave the active selection on G5, as in the saved workbook
$ws.Range("G5").Select()
